$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - the localization status report was
# regenerated: items that were previously "In Translation" are now
# "Ready for handoff", and the HO Xliff generation timestamps were bumped to
# reflect the new report run. Columns whose text grew wider were re-sized to
# fit the new content.
# ---------------------------------------------------------------------------

# --- Overview sheet ---------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 09:21:37"

# Columns E (zh-cn) and F (de-de) need to widen to fit "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 09:21:32"

# Column C (Status) widens to fit "Ready for handoff".
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet -------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 09:21:37"

# Column C (Status) widens to fit "Ready for handoff".
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
